$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.153878919646918
$ws.Range("C2").Value = 1.42188496435751
$ws.Range("B3").Value = 8.572081016669772
$ws.Range("C3").Value = 2.830981203228037
$ws.Range("B4").Value = 9.200494610421321
$ws.Range("C4").Value = 4.446453976652502
$ws.Range("B5").Value = 10.50103078494874
$ws.Range("C5").Value = 5.832339325672852
$ws.Range("B6").Value = 18.967763382494
$ws.Range("C6").Value = 7.151241821523191
$ws.Range("B7").Value = 19.23972948902364
$ws.Range("C7").Value = 8.633334693791143
$ws.Range("B8").Value = 25.69875723445718
$ws.Range("C8").Value = 9.887305346799046
$ws.Range("B9").Value = 27.71724102012856
$ws.Range("C9").Value = 11.09873488619389
$ws.Range("B10").Value = 28.91682808155998
$ws.Range("C10").Value = 12.49507968074602
$ws.Range("B11").Value = 29.32629327005424
$ws.Range("C11").Value = 13.95047353686999
$ws.Range("B12").Value = 33.6071320299563
$ws.Range("C12").Value = 15.4039481908298
$ws.Range("B13").Value = 33.73200775847001
$ws.Range("C13").Value = 16.68865640904783
$ws.Range("B14").Value = 34.26943263270257
$ws.Range("C14").Value = 17.99496937955135
$ws.Range("B15").Value = 36.53323051239772
$ws.Range("C15").Value = 19.37827436681111
$ws.Range("B16").Value = 37.89174466619554
$ws.Range("C16").Value = 20.63830136446376
$ws.Range("B17").Value = 38.35016081964368
$ws.Range("C17").Value = 21.89283980788699
$ws.Range("B18").Value = 40.13946161702658
$ws.Range("C18").Value = 23.43506501934904
$ws.Range("B19").Value = 43.45331736390975
$ws.Range("C19").Value = 24.78048935579151
$ws.Range("B20").Value = 51.40063609145808
$ws.Range("C20").Value = 26.1791522549488
$ws.Range("B21").Value = 56.37645327532704
$ws.Range("C21").Value = 27.45865051743643
$ws.Range("B22").Value = 58.43197744858656
$ws.Range("C22").Value = 29.39983088374621
$ws.Range("B23").Value = 61.56425521650574
$ws.Range("C23").Value = 30.73557923526711
$ws.Range("B24").Value = 63.89383001009763
$ws.Range("C24").Value = 31.88463036038733
$ws.Range("B25").Value = 64.55994840438292
$ws.Range("C25").Value = 33.20342347340743
$ws.Range("B26").Value = 65.4669069213049
$ws.Range("C26").Value = 34.39435739783485
$ws.Range("B27").Value = 69.08919683735955
$ws.Range("C27").Value = 35.78777588225208
$ws.Range("B28").Value = 70.14177588856035
$ws.Range("C28").Value = 36.9973777818394
$ws.Range("B29").Value = 71.04152908506028
$ws.Range("C29").Value = 38.25409108435428
$ws.Range("B30").Value = 71.94272657987695
$ws.Range("C30").Value = 39.35656368941497
$ws.Range("B31").Value = 74.16414854777514
$ws.Range("C31").Value = 40.60033102221058
$ws.Range("B32").Value = 75.117084458447
$ws.Range("C32").Value = 41.7705375026037
$ws.Range("B33").Value = 77.24657068392382
$ws.Range("C33").Value = 43.13818193468219
$ws.Range("B34").Value = 81.55808260037792
$ws.Range("C34").Value = 44.45433001447621
$ws.Range("B35").Value = 83.00676797187819
$ws.Range("C35").Value = 45.70323077914892
$ws.Range("B36").Value = 85.10352364378639
$ws.Range("C36").Value = 47.01634630455409
$ws.Range("B37").Value = 85.24309963097372
$ws.Range("C37").Value = 48.16719134883965
$ws.Range("B38").Value = 86.23240935339048
$ws.Range("C38").Value = 49.60940669504738
$ws.Range("B39").Value = 88.4642341938564
$ws.Range("C39").Value = 51.20928921728122
$ws.Range("B40").Value = 88.82085648802716
$ws.Range("C40").Value = 52.63974194663394
$ws.Range("B41").Value = 89.00935976293485
$ws.Range("C41").Value = 53.87739762733833
$ws.Range("B42").Value = 89.82541420659417
$ws.Range("C42").Value = 55.28760591085513
$ws.Range("B43").Value = 91.42222524487751
$ws.Range("C43").Value = 56.62860274626167
$ws.Range("B44").Value = 91.91879571052451
$ws.Range("C44").Value = 57.99631718147283
$ws.Range("B45").Value = 95.39337758099614
$ws.Range("C45").Value = 59.57284302940161
$ws.Range("B46").Value = 95.50208326931305
$ws.Range("C46").Value = 61.0361062415381
$ws.Range("B47").Value = 96.095439941102
$ws.Range("C47").Value = 62.30560557392884
$ws.Range("B48").Value = 97.52877653404371
$ws.Range("C48").Value = 63.8651357521776
$ws.Range("B49").Value = 97.64767765213379
$ws.Range("C49").Value = 65.10867722405796
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 99.58789267194999
$ws.Range("C50").Value = 66.34743703584
